# Apply "hybrid bold + color" highlighting to quantitative metrics
# (percentages, dollar amounts, large numbers) across the resume.
#
# For each target paragraph (1-based Paragraphs index) we split the
# plain-text run into multiple runs: the metric substrings get
# Bold = True and Font.Color = RGB(0x2C,0x3E,0x50); everything else
# is left as plain text.

$d = $word.ActiveDocument

$HighlightColor = 5258796   # RGB(0x2C, 0x3E, 0x50) == 0x00503E2C (BGR-packed OLE color)

function Highlight-Metric {
    param($ParaIndex, [string]$Metric)

    $para = $d.Paragraphs.Item($ParaIndex)
    $searchRng = $para.Range.Duplicate
    $found = $searchRng.Find.Execute($Metric, $false, $false, $false, $false, $false, $true, 0, $false, "", 0)
    if ($found) {
        $searchRng.Font.Bold = 1
        $searchRng.Font.Color = $HighlightColor
    }
}

# Paragraph 10: "... classification accuracy from 23% to 64%"
Highlight-Metric 10 "23%"
Highlight-Metric 10 "64%"

# Paragraph 12: "... margin of error from ±4.2% to ±2.1%, increasing voter
# turnout prediction accuracy from 71% to 87%, ..."
Highlight-Metric 12 "±4.2%"
Highlight-Metric 12 "±2.1%"
Highlight-Metric 12 "71%"
Highlight-Metric 12 "87%"

# Paragraph 13: "... reduced mapping costs by 73.5%, saving campaigns and
# organizations $4.7M and enabling ..."
Highlight-Metric 13 "73.5%"
Highlight-Metric 13 "$4.7M"

# Paragraph 14: "... sub-economy valued over $2 trillion"
Highlight-Metric 14 "$2"

# Paragraph 20: "... reducing processing time by 57%"
Highlight-Metric 20 "57%"

# Paragraph 85: "• 178% accuracy improvement in racial classification algorithms"
Highlight-Metric 85 "178%"

# Paragraph 86: "... reducing mapping costs 73.5%"
Highlight-Metric 86 "73.5%"

# Paragraph 87: "• $4.7M savings enabled nonprofit access"
Highlight-Metric 87 "$4.7M"

# Paragraph 88: "... system serving 12,847 analysts across 89 organizations"
Highlight-Metric 88 "12,847"

# Paragraph 90: "... margin of error from ±4.2% to ±2.1%"
Highlight-Metric 90 "±4.2%"
Highlight-Metric 90 "±2.1%"

# Paragraph 91: "• Increased voter turnout prediction accuracy from 71% to 87%"
Highlight-Metric 91 "71%"
Highlight-Metric 91 "87%"
